$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-15 (columns B and E change)
$ws.Range("B2").Value = 21469
$ws.Range("E2").Value = 15259

$ws.Range("B3").Value = 20621
$ws.Range("E3").Value = 14648

$ws.Range("B4").Value = 19828
$ws.Range("E4").Value = 13520

$ws.Range("B5").Value = 19009
$ws.Range("E5").Value = 12691

$ws.Range("B6").Value = 18488
$ws.Range("E6").Value = 12274

$ws.Range("B7").Value = 18753
$ws.Range("E7").Value = 12200

$ws.Range("B8").Value = 18740
$ws.Range("E8").Value = 12590

$ws.Range("B9").Value = 19624
$ws.Range("E9").Value = 12971

$ws.Range("B10").Value = 22373
$ws.Range("E10").Value = 13835

$ws.Range("B11").Value = 24582
$ws.Range("E11").Value = 14487

$ws.Range("B12").Value = 25034
$ws.Range("E12").Value = 14850

$ws.Range("B13").Value = 24749
$ws.Range("E13").Value = 14469

$ws.Range("B14").Value = 24350
$ws.Range("E14").Value = 14333

$ws.Range("B15").Value = 25698
$ws.Range("E15").Value = 15180

# Fill in new rows 16-19 with data (previously empty placeholder rows)
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 26016
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 15492
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 25623
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 15143
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 24213
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 14781
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 21283
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 14963
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
